$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.244.78"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.869.20"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'0.7103"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "'241.37"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.3100"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "'0.07682"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "'24.94"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'0.08369"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "1.879.47"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "'5.201"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "'0.7088"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "'91.06"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "29.249.10"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'0.000008253"
$ws.Range("E17").Value = "  +5.60%  "
$ws.Range("D18").Value = "'5.926"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "'242.22"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").Value = "2.130.41"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'13.14"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'7.809"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'0.1628"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").Value = "'162.99"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "'8.991"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "'4.399"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'4.318"
$ws.Range("E31").Value = "  +5.12%  "
$ws.Range("D32").Value = "'1.280"
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").Value = "'0.05239"
$ws.Range("D34").Value = "'1.919"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'0.7501"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").Value = "'1.169"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").Value = "'2.684"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.01855"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'2.714"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "1.153.87"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "'6.353"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("D42").Value = "'72.92"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "'0.8847"
$ws.Range("D44").Value = "'104.27"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D49").Value = "'9.344"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'0.4285"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'0.9980"
$ws.Range("E51").Value = "  +0.02%  "
